# Finish adding the standard metadata columns (property_category, category,
# date, legislator_name, legislator_id, source_file, index) to the "債務"
# (debt) and "事業投資" (business investment) sheets, matching the layout
# already used on the other sheets of this workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 6: 債務 (Debt)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)

# Header row (row 1): species, debtor, owner, total, register_date,
# register_reason, then the shared metadata headers.
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Row 2 (index 110)
$ws.Range("B2").Value = "股票質押借款"
$ws.Range("C2").Value = "蔡慧敏"
$ws.Range("D2").Value = "第一商業銀行北桃桃圜縣桃園市三民路"
$ws.Range("E2").Value = 50000000
$ws.Range("F2").Value = "98年06月01日"
$ws.Range("G2").Value = "借入現金"
$ws.Range("H2").Value = "debt"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2012-02-29"
$ws.Range("K2").Value = "陳根德"
$ws.Range("L2").Value = 833
$ws.Range("M2").Value = "tmp28cf1"
$ws.Range("N2").Value = 110

# Row 3 (index 111)
$ws.Range("B3").Value = "股票質押借款"
$ws.Range("C3").Value = "蔡慧敏"
$ws.Range("D3").Value = "聯邦商業銀行桃圜桃園縣桃圜市中山路"
$ws.Range("E3").Value = 20000000
$ws.Range("F3").Value = "100年5月"
$ws.Range("G3").Value = "借人現金"
$ws.Range("H3").Value = "debt"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").Value = "2012-02-29"
$ws.Range("K3").Value = "陳根德"
$ws.Range("L3").Value = 833
$ws.Range("M3").Value = "tmp28cf1"
$ws.Range("N3").Value = 111

# Row 4 (index 112)
$ws.Range("B4").Value = "房屋貸款"
$ws.Range("C4").Value = "陳根德"
$ws.Range("D4").Value = "合作金庫商業銀行桃圜桃園縣桃圜市中正路"
$ws.Range("E4").Value = 22366025
$ws.Range("F4").Value = "93年01月12日"
$ws.Range("G4").Value = "房貸"
$ws.Range("H4").Value = "debt"
$ws.Range("I4").Value = "normal"
$ws.Range("J4").Value = "2012-02-29"
$ws.Range("K4").Value = "陳根德"
$ws.Range("L4").Value = 833
$ws.Range("M4").Value = "tmp28cf1"
$ws.Range("N4").Value = 112

$ws.Range("A1:N4").EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 7: 事業投資 (Business investment)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(7)

# Header row (row 1): owner, company, address, total, register_date,
# register_reason, then the shared metadata headers.
$ws2.Range("B1").Value = "owner"
$ws2.Range("C1").Value = "company"
$ws2.Range("D1").Value = "address"
$ws2.Range("E1").Value = "total"
$ws2.Range("F1").Value = "register_date"
$ws2.Range("G1").Value = "register_reason"
$ws2.Range("H1").Value = "property_category"
$ws2.Range("I1").Value = "category"
$ws2.Range("J1").Value = "date"
$ws2.Range("K1").Value = "legislator_name"
$ws2.Range("L1").Value = "legislator_id"
$ws2.Range("M1").Value = "source_file"
$ws2.Range("N1").Value = "index"

# Row 2 (index 117)
$ws2.Range("B2").Value = "陳根德"
$ws2.Range("C2").Value = "坤和雷射品切割(股）"
$ws2.Range("D2").Value = "桃園縣龜山鄉舊路村9鄰振興路998號"
$ws2.Range("E2").Value = 2000000
$ws2.Range("F2").Value = "94年03月16日"
$ws2.Range("G2").Value = "投資"
$ws2.Range("H2").Value = "investment"
$ws2.Range("I2").Value = "normal"
$ws2.Range("J2").Value = "2012-02-29"
$ws2.Range("K2").Value = "陳根德"
$ws2.Range("L2").Value = 833
$ws2.Range("M2").Value = "tmp28cf1"
$ws2.Range("N2").Value = 117

# Row 3 (index 118)
$ws2.Range("B3").Value = "陳根德"
$ws2.Range("C3").Value = "建道營造(股）"
$ws2.Range("D3").Value = "桃園縣桃園市四維街5號"
$ws2.Range("E3").Value = 45000000
$ws2.Range("F3").Value = "91年08月05日"
$ws2.Range("G3").Value = "投資"
$ws2.Range("H3").Value = "investment"
$ws2.Range("I3").Value = "normal"
$ws2.Range("J3").Value = "2012-02-29"
$ws2.Range("K3").Value = "陳根德"
$ws2.Range("L3").Value = 833
$ws2.Range("M3").Value = "tmp28cf1"
$ws2.Range("N3").Value = 118

# Row 4 (index 119)
$ws2.Range("B4").Value = "蔡慧敏"
$ws2.Range("C4").Value = "建道營造(股）"
$ws2.Range("D4").Value = "桃圜縣桃園市四維街5號"
$ws2.Range("E4").Value = 41700000
$ws2.Range("F4").Value = "91年08月05日"
$ws2.Range("G4").Value = "投資"
$ws2.Range("H4").Value = "investment"
$ws2.Range("I4").Value = "normal"
$ws2.Range("J4").Value = "2012-02-29"
$ws2.Range("K4").Value = "陳根德"
$ws2.Range("L4").Value = 833
$ws2.Range("M4").Value = "tmp28cf1"
$ws2.Range("N4").Value = 119

# Row 5 (index 120)
$ws2.Range("B5").Value = "蔡慧敏"
$ws2.Range("C5").Value = "天天美企業(有）"
$ws2.Range("D5").Value = "桃圜縣桃園市四維街5號"
$ws2.Range("E5").Value = 500000
$ws2.Range("F5").Value = "92年04月15曰"
$ws2.Range("G5").Value = "投資"
$ws2.Range("H5").Value = "investment"
$ws2.Range("I5").Value = "normal"
$ws2.Range("J5").Value = "2012-02-29"
$ws2.Range("K5").Value = "陳根德"
$ws2.Range("L5").Value = 833
$ws2.Range("M5").Value = "tmp28cf1"
$ws2.Range("N5").Value = 120

# Row 6 (index 122)
$ws2.Range("B6").Value = "蔡慧敏"
$ws2.Range("C6").Value = "育嘉建設(有）"
$ws2.Range("D6").Value = "桃圜縣桃園市四维街5號"
$ws2.Range("E6").Value = 1000000
$ws2.Range("F6").Value = "95年12月20日"
$ws2.Range("G6").Value = "投資"
$ws2.Range("H6").Value = "investment"
$ws2.Range("I6").Value = "normal"
$ws2.Range("J6").Value = "2012-02-29"
$ws2.Range("K6").Value = "陳根德"
$ws2.Range("L6").Value = 833
$ws2.Range("M6").Value = "tmp28cf1"
$ws2.Range("N6").Value = 122

$ws2.Range("A1:N6").EntireColumn.AutoFit() | Out-Null
